$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(4, 6).Value = -3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(11, 6).Value = -5
$ws.Cells.Item(13, 6).Value = -2
$ws.Cells.Item(14, 6).Value = -3
$ws.Cells.Item(17, 6).Value = -5
$ws.Cells.Item(18, 6).Value = -5
$ws.Cells.Item(22, 6).Value = -5
$ws.Cells.Item(32, 6).Value = -1
$ws.Cells.Item(43, 6).Value = 3
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(46, 6).Value = -1
$ws.Cells.Item(47, 6).Value = 1
$ws.Cells.Item(51, 6).Value = 10
$ws.Cells.Item(52, 6).Value = 4
$ws.Cells.Item(54, 6).Value = -4
$ws.Cells.Item(55, 6).Value = 5
$ws.Cells.Item(56, 6).Value = 11
$ws.Cells.Item(57, 6).Value = -2
$ws.Cells.Item(58, 6).Value = 8
$ws.Cells.Item(63, 6).Value = -5
$ws.Cells.Item(65, 6).Value = -1
$ws.Cells.Item(67, 6).Value = -2
